# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-row price records (columns D, L-T)
# across the existing rows 2-11, while columns A-C and E-K (market/product
# identity columns) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that get reshuffled, keyed by
# original row number, so we can redistribute them to their new rows.
$cols = @("D","L","M","N","O","P","Q","R","S","T")

$snapshot = @{}
for ($r = 2; $r -le 11; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $row
}

# Mapping of new row -> source row whose D/L-T values it should now hold.
$mapping = @{
    2  = 7
    3  = 8
    4  = 2
    5  = 6
    6  = 9
    7  = 10
    8  = 5
    9  = 11
    10 = 3
    11 = 4
}

foreach ($newRow in 2..11) {
    $srcRow = $mapping[$newRow]
    $srcData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $srcData[$c]
    }
}
